$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated KPI results (B2:D6) per the case-study re-run
$ws.Range("B2").Value = 84666.358663243998
$ws.Range("C2").Value = 61150.869842430897
$ws.Range("D2").Value = 342.99568113491699

$ws.Range("B3").Value = 82856.030387143095
$ws.Range("C3").Value = 59413.831379038602
$ws.Range("D3").Value = 416.28549384412503

$ws.Range("B4").Value = 80743.654556177498
$ws.Range("C4").Value = 57316.646779459297
$ws.Range("D4").Value = 431.47672523089102

$ws.Range("B5").Value = 76958.598795020604
$ws.Range("C5").Value = 53635.017615921402
$ws.Range("D5").Value = 534.90332284982401

$ws.Range("B6").Value = 70847.511616987496
$ws.Range("C6").Value = 47786.348328721302
$ws.Range("D6").Value = 797.32121368105697

# Update the active selection to match the saved view state
$ws.Range("C13").Select()
